$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS")

# Helper cell used to build date-look-alike text via a formula result
# (so it never goes through Excel's "looks like a date" literal-entry
# autodetection), then copied as a value onto the real destination cell.
$helper = $ws.Cells.Item(500, 1)

function Set-TextValue($cell, [string]$text) {
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# --- Row 10: existing row - tweak B10's time value slightly, and bring
#     the rest of the row onto the column's default style (s=5) by
#     clearing + re-writing the content. ---
$ws.Cells.Item(10, 1).ClearContents()
Set-TextValue $ws.Cells.Item(10, 1) "2021-06-11"

$ws.Cells.Item(10, 2).Value = 44358.65313825232

$ws.Cells.Item(10, 3).ClearContents()
Set-TextValue $ws.Cells.Item(10, 3) "pavan_demo_145"

$ws.Cells.Item(10, 4).ClearContents()
$ws.Cells.Item(10, 4).Value = 96

$ws.Cells.Item(10, 5).ClearContents()
$ws.Cells.Item(10, 5).Value = 95

$ws.Cells.Item(10, 6).ClearContents()
$ws.Cells.Item(10, 6).Value = 1

$ws.Cells.Item(10, 7).ClearContents()
$ws.Cells.Item(10, 7).Value = 2.51

# --- Row 11: previously a blank placeholder row; fill with the new
#     interview-run data. ---
Set-TextValue $ws.Cells.Item(11, 1) "2021-06-16"

$ws.Cells.Item(11, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 2).Value = 44363.68435071759

Set-TextValue $ws.Cells.Item(11, 3) "live_145_hotfix"

$ws.Cells.Item(11, 4).Value = 96
$ws.Cells.Item(11, 5).Value = 95
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.46

# --- Row 12: previously a blank placeholder row; fill with the new
#     interview-run data (hotfix 2). ---
Set-TextValue $ws.Cells.Item(12, 1) "2021-06-16"

$ws.Cells.Item(12, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 2).Value = 44363.80313471368

Set-TextValue $ws.Cells.Item(12, 3) "live_145_hf2"

$ws.Cells.Item(12, 4).Value = 96
$ws.Cells.Item(12, 5).Value = 94
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(12, 7).Value = 2.49

# Clean up helper cell
$helper.Clear()
